# Initial deployment of Azure VMs from Excel
# Adds a second VM row (row 3) to the "VM Details" sheet, duplicating the
# values/format of row 2 except for Resource Group, VM Name and the
# network interface name which are unique to the new VM.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 values (mirrors row 2, with 3 new unique values) ---
# NOTE: reading a cell's value back out requires calling Value() as a
# method (with parens) in this COM host - a bare `.Value` access returns
# the property descriptor rather than invoking the getter.
$ws.Range("A3").Value = $ws.Range("A2").Value()
$ws.Range("B3").Value = "jf-rg-001"
$ws.Range("C3").Value = "jf-vm-001"
$ws.Range("D3").Value = $ws.Range("D2").Value()
$ws.Range("F3").Value = $ws.Range("F2").Value()
$ws.Range("G3").Value = $ws.Range("G2").Value()
$ws.Range("H3").Value = $ws.Range("H2").Value()
$ws.Range("I3").Value = $ws.Range("I2").Value()
$ws.Range("J3").Value = $ws.Range("J2").Value()
$ws.Range("K3").Value = $ws.Range("K2").Value()
$ws.Range("L3").Value = $ws.Range("L2").Value()
$ws.Range("M3").Value = $ws.Range("M2").Value()
$ws.Range("N3").Value = $ws.Range("N2").Value()
$ws.Range("O3").Value = $ws.Range("O2").Value()
$ws.Range("P3").Value = "jf-vm-001754"
$ws.Range("W3").Value = $ws.Range("W2").Value()
$ws.Range("X3").Value = $ws.Range("X2").Value()

# --- Match row 2 formatting (fonts/styles) for the styled cells ---
$ws.Range("I2:J2").Copy()
$ws.Range("I3:J3").PasteSpecial(-4122)

$ws.Range("N2").Copy()
$ws.Range("N3").PasteSpecial(-4122)

$ws.Range("L2").Copy()
$ws.Range("L3").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Hyperlink the new password cell, same as L2 ---
$ws.Hyperlinks.Add($ws.Range("L3"), "mailto:P@ssword1234") | Out-Null

# re-apply L2's cell style, since adding the hyperlink restyles the cell
$ws.Range("L2").Copy()
$ws.Range("L3").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Row height to match row 2 ---
$ws.Rows(3).RowHeight = $ws.Rows(2).RowHeight

# --- Update the view/selection to the new row ---
$ws.Range("AH3").Select() | Out-Null
$excel.ActiveWindow.ScrollColumn = 21
$excel.ActiveWindow.ScrollRow = 1

Write-Host "Row 3 added"
